# Auto-generated edit script: update cryptos.xlsx D/E (and a few B/C) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell {
    param($range, [string]$val)
    # Force the cell to store the literal text, preserving formatting
    # such as trailing zeros / multiple dots, without leaving a permanent
    # number-format override behind once done.
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "66.231.79"
$ws.Range("E2").Value = "  -3.45%  "
Set-TextCell $ws.Range("D3") "3.777.39"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextCell $ws.Range("D5") "417.27"
$ws.Range("E5").Value = "  -0.66%  "
Set-TextCell $ws.Range("D6") "126.90"
$ws.Range("E6").Value = "  -7.62%  "
Set-TextCell $ws.Range("D7") "3.774.86"
$ws.Range("E7").Value = "  +0.00%  "
Set-TextCell $ws.Range("D8") "0.598"
$ws.Range("E8").Value = "  -7.96%  "
Set-TextCell $ws.Range("D9") "0.999"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -8.47%  "
$ws.Range("E11").Value = "  -15.45%  "
$ws.Range("E12").Value = "  -17.06%  "
Set-TextCell $ws.Range("D13") "39.53"
$ws.Range("E13").Value = "  -8.49%  "
Set-TextCell $ws.Range("D14") "4.428.49"
$ws.Range("E14").Value = "  +1.04%  "
Set-TextCell $ws.Range("D15") "16.12"
$ws.Range("E15").Value = "  +22.89%  "
Set-TextCell $ws.Range("D16") "9.91"
$ws.Range("E16").Value = "  -5.70%  "
$ws.Range("E17").Value = "  -1.80%  "
Set-TextCell $ws.Range("D18") "3.811.53"
$ws.Range("E18").Value = "  -0.07%  "
Set-TextCell $ws.Range("D19") "19.24"
$ws.Range("E19").Value = "  -6.37%  "
Set-TextCell $ws.Range("D20") "66.524.27"
$ws.Range("E20").Value = "  -2.71%  "
$ws.Range("E21").Value = "  -6.06%  "
Set-TextCell $ws.Range("D22") "402.60"
$ws.Range("E22").Value = "  -9.71%  "
Set-TextCell $ws.Range("D23") "14.03"
$ws.Range("E23").Value = "  -9.49%  "
Set-TextCell $ws.Range("D24") "83.00"
$ws.Range("E24").Value = "  -8.40%  "
$ws.Range("E25").Value = "  -4.65%  "
Set-TextCell $ws.Range("D26") "5.73"
$ws.Range("E26").Value = "  +12.27%  "
Set-TextCell $ws.Range("D27") "36.35"
$ws.Range("E27").Value = "  -5.38%  "
Set-TextCell $ws.Range("D28") "3.14"
$ws.Range("E28").Value = "  -5.56%  "
Set-TextCell $ws.Range("D29") "9.22"
$ws.Range("E29").Value = "  -8.74%  "
Set-TextCell $ws.Range("D30") "705.46"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws.Range("D31") "7.62"
$ws.Range("E31").Value = "  +6.43%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws.Range("D32") "2.75"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("E35").Value = "  -9.65%  "
Set-TextCell $ws.Range("D36") "0.998"
$ws.Range("E36").Value = "  -0.19%  "
Set-TextCell $ws.Range("D37") "37.14"
$ws.Range("E37").Value = "  -10.21%  "
Set-TextCell $ws.Range("D38") "54.56"
$ws.Range("E38").Value = "  -6.05%  "
Set-TextCell $ws.Range("D39") "0.0₃0754"
$ws.Range("E39").Value = "  +5.98%  "
Set-TextCell $ws.Range("D40") "0.0446"
$ws.Range("E40").Value = "  -8.84%  "
$ws.Range("E41").Value = "  -4.16%  "
Set-TextCell $ws.Range("D42") "1.01"
$ws.Range("E42").Value = "  +1.07%  "
Set-TextCell $ws.Range("D43") "4.51"
$ws.Range("E43").Value = "  +4.32%  "
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("E45").Value = "  -10.18%  "
Set-TextCell $ws.Range("D46") "3.27"
$ws.Range("E46").Value = "  -3.26%  "
Set-TextCell $ws.Range("D47") "143.50"
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("E48").Value = "  -3.24%  "
Set-TextCell $ws.Range("D49") "25.50"
$ws.Range("E49").Value = "  -6.90%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell $ws.Range("D50") "4.568.86"
$ws.Range("E50").Value = "  +10.01%  "
$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws.Range("D51") "2.50"
$ws.Range("E51").Value = "  -4.29%  "
